$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column B (rows 2-89) holds date-like text such as "2012-01-01".
# Reformat each of those values from ISO "YYYY-01-01" to US "01/01/YYYY",
# keeping them as plain text (not real Excel dates).
$firstRow = 2
$lastRow = 89

$dateRange = $ws.Range("B2:B89")

# Mark the range as text first so that writing a "01/01/2012"-shaped string
# into it isn't auto-recognized and converted into a date serial number.
$dateRange.NumberFormat = "@"

for ($row = $firstRow; $row -le $lastRow; $row++) {
    $cell = $ws.Cells.Item($row, 2)
    $old = $cell.Value2

    if ($old -match "^(\d{4})-01-01$") {
        $year = $matches[1]
        $cell.Value = "01/01/" + $year
    }
}

# Drop the temporary text-number-format again so the cells end up with no
# explicit style, exactly like before the edit.
$dateRange.ClearFormats()
